# "combine input and output sankey"
# expenses (sheet1) gets two wider helper columns (B, C) for the combined
# sankey links; income (sheet2) gets recoded so its "Job"/"Erbe" rows feed
# the same link structure ("Novatec"/"gehalt", "a"/"b").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # expenses
$ws2 = $wb.Worksheets.Item(2)   # income

# --- income (sheet2): recode row 2 + row 3 ---------------------------------
$ws2.Range("B2").Value = "Novatec"
$ws2.Range("C2").Value = "gehalt"

$ws2.Range("A3").ClearContents()
$ws2.Range("B3").Value = "a"
$ws2.Range("C3").Value = "b"

# --- expenses (sheet1): widen the helper columns used by the combined chart -
$ws1.Columns.Item(2).ColumnWidth = 19.666666666666668   # -> stored width ~20.5546875
$ws1.Columns.Item(3).ColumnWidth = 29.666666666666668   # -> stored width ~30.5546875

# --- selections / active sheet ---------------------------------------------
# expenses: select D2:D7 (no longer the active/tabSelected sheet)
$ws1.Range("D2:D7").Select()

# income becomes the active (tabSelected) sheet, selection moves to M18
$ws2.Activate()
$ws2.Range("M18").Select()
